$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: update title and link
$ws.Range("D28").Value = "Let's do MuJoCo"
$ws.Range("E28").Value = "https://ropiens.tistory.com/168"

# Row 37: update title and link
$ws.Range("D37").Value = "[Paper Review] SimCSE : Simple Contrastive Learning of Sentence Embeddings"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1875&mod=document&pageid=1"

# Row 50: update title and link
$ws.Range("D50").Value = "수학과 기계학습"
$ws.Range("E50").Value = "http://incredible.egloos.com/7531273"
